$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Section 1: Poly kernel ROC points (rows 5-10), chart1 data ---
$ws.Range("E5").Value = 0.478
$ws.Range("F5").Value = 0.176

$ws.Range("E6").Value = 0.478
$ws.Range("F6").Value = 0.176

$ws.Range("E7").Value = 0.478
$ws.Range("F7").Value = 0.176

$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 1

$ws.Range("E9").Value = 0.348
$ws.Range("F9").Value = 0.029

$ws.Range("E10").Value = 1
$ws.Range("F10").Value = 1

# --- Section 2: RBF kernel ROC points (rows 17-34), chart2 data ---
$ws.Range("C17").Value = "RBF"

$ws.Range("E17").Value = 1
$ws.Range("F17").Value = 1

$ws.Range("E18").Value = 1
$ws.Range("F18").Value = 1

$ws.Range("E19").Value = 1
$ws.Range("F19").Value = 1

$ws.Range("E20").Value = 1
$ws.Range("F20").Value = 1

$ws.Range("E21").Value = 1
$ws.Range("F21").Value = 1

$ws.Range("E22").Value = 1
$ws.Range("F22").Value = 1

$ws.Range("E23").Value = 1
$ws.Range("F23").Value = 1

$ws.Range("E24").Value = 1
$ws.Range("F24").Value = 1

$ws.Range("E25").Value = 1
$ws.Range("F25").Value = 1

$ws.Range("E26").Value = 1
$ws.Range("F26").Value = 1

$ws.Range("E27").Value = 1
$ws.Range("F27").Value = 1

$ws.Range("E28").Value = 1
$ws.Range("F28").Value = 1

$ws.Range("E29").Value = 1
$ws.Range("F29").Value = 1

$ws.Range("E30").Value = 0.304
$ws.Range("F30").Value = 0.029

$ws.Range("E31").Value = 0.348
$ws.Range("F31").Value = 0.029

$ws.Range("E32").Value = 0.435
$ws.Range("F32").Value = 0.059

$ws.Range("E33").Value = 0.522
$ws.Range("F33").Value = 0.118

$ws.Range("E34").Value = 0.478
$ws.Range("F34").Value = 0.118

# --- Section 3: selection cell matches author's final cursor position ---
$ws.Range("C18").Select() | Out-Null
